# Commit: "Electrolyser Data Updated(name of the sheet) and excel reader
# adapted for electrolyser. Now 2 libraries arrives at technical comparison"
#
# The workbook's second tab ("Sheet2") - which holds the electrolyser
# technical-comparison data used by the external reader script - is
# renamed to "Script (Main)" to reflect its new role as the main script
# sheet. The active selection on that sheet is also moved to D14.

$wb = $excel.ActiveWorkbook

# Locate the worksheet named "Sheet2" (fall back to the second sheet if the
# workbook's tab order/naming differs) and rename it.
$sheet2 = $null
try {
    $sheet2 = $wb.Worksheets.Item("Sheet2")
} catch {
    $sheet2 = $null
}
if (-not $sheet2) {
    $sheet2 = $wb.Worksheets.Item(2)
}

$sheet2.Name = "Script (Main)"

# Make it the active tab and move the selection to D14, matching the
# author's final cursor position after the edit.
$sheet2.Activate()
$sheet2.Range("D14").Select()
